$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3345
$ws1.Range("F5").Value = 2418
$ws1.Range("F7").Value = 338
$ws1.Range("F8").Value = 1369
$ws1.Range("F14").Value = 95
$ws1.Range("F16").Value = 8444
$ws1.Range("F18").Value = 2481
$ws1.Range("F27").Value = 1978
$ws1.Range("F32").Value = 1913
$ws1.Range("F43").Value = 110

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3345
$ws4.Range("F5").Value = 2418
$ws4.Range("F7").Value = 338
$ws4.Range("F8").Value = 1369
$ws4.Range("F14").Value = 95
$ws4.Range("F16").Value = 8444
$ws4.Range("F18").Value = 2481
$ws4.Range("F28").Value = 1978
$ws4.Range("F32").Value = 1913
$ws4.Range("F47").Value = 110
